$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.519.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.020.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.09%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.020.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.19%  "
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("E13").Value = "  +3.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.520.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.513.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.020.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "448.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.690"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.52%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.11%  "
$ws.Range("E32").Value = "  +4.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("E35").Value = "  +10.41%  "
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.07%  "
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("E42").Value = "  +4.59%  "
$ws.Range("E43").Value = "  +9.20%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.76%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "390.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0355"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.737.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("E51").Value = "  +0.12%  "